$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.898.51"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "2.165.75"
$ws.Range("E3").Value = "  -2.92%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.07"
$ws.Range("E5").Value = "  -2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.97"
$ws.Range("E7").Value = "  -7.88%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.97"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0921"
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "35.41"
$ws.Range("E12").Value = "  -16.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "2.486.30"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.853"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.17"
$ws.Range("E17").Value = "  -5.54%  "
$ws.Range("D18").Value = "2.169.80"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").Value = "40.735.76"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.07"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "228.11"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  -8.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.56"
$ws.Range("E25").Value = "  +12.98%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.68"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E28").Value = "  -5.44%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.87"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  -9.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.09"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.61"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0734"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "25.43"
$ws.Range("E36").Value = "  -4.41%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0297"
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.16"
$ws.Range("E40").Value = "  -5.51%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.60"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.43"
$ws.Range("E42").Value = "  -9.99%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.81"
$ws.Range("E43").Value = "  -5.19%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.63"
$ws.Range("E44").Value = "  -16.21%  "
$ws.Range("E45").Value = "  -12.44%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.41"
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0983"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.14"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("E51").Value = "  -1.35%  "
